# Fruta / hortaliza, semanal
# New weekly price observations (date 44578) are inserted for this market,
# pushing the two previously-most-recent rows (and everything below them)
# down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 48-49, shifting existing rows 48:51 down to 50:53.
$ws.Rows("48:49").Insert()

# --- Row 48: new observation, Calameño / Primera ---
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44578
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112027
$ws.Cells.Item(48, 7).Value = "Melón"
$ws.Cells.Item(48, 8).Value = "Calameño"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 100
$ws.Cells.Item(48, 11).Value = 8000
$ws.Cells.Item(48, 12).Value = 9000
$ws.Cells.Item(48, 13).Value = 8500
$ws.Cells.Item(48, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 708
$ws.Cells.Item(48, 17).Value = 12
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# --- Row 49: new observation, Tuna / Segunda ---
$ws.Cells.Item(49, 1).Value = 1
$ws.Cells.Item(49, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value = 44578
$ws.Cells.Item(49, 5).Value = 15
$ws.Cells.Item(49, 6).Value = 100112027
$ws.Cells.Item(49, 7).Value = "Melón"
$ws.Cells.Item(49, 8).Value = "Tuna"
$ws.Cells.Item(49, 9).Value = "Segunda"
$ws.Cells.Item(49, 10).Value = 80
$ws.Cells.Item(49, 11).Value = 6000
$ws.Cells.Item(49, 12).Value = 7000
$ws.Cells.Item(49, 13).Value = 6500
$ws.Cells.Item(49, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 271
$ws.Cells.Item(49, 17).Value = 24
$ws.Cells.Item(49, 18).Value = "Hortaliza"
